$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCD")

# Insert two new columns before column D (the data shifts from D:K to F:M,
# matching how the workbook now carries two more quarters of history).
$ws.Range("D1:E1").EntireColumn.Insert()

# Carry over number formatting/styles for the freshly inserted D:E columns
# from the (now shifted) old D:E columns, which live in F:G after the insert.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the column widths of the newly inserted columns to their neighbors.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Populate the two new quarters of data (new columns D and E).
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 5163000
$ws.Range("E8").Value = 5369400
$ws.Range("D9").Value = 2466300
$ws.Range("E9").Value = 2547300
$ws.Range("D10").Value = 2696700
$ws.Range("E10").Value = 2822100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 90500
$ws.Range("E14").Value = -67500
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 3163500
$ws.Range("E17").Value = 2951700
$ws.Range("D18").Value = 1999500
$ws.Range("E18").Value = 2417700
$ws.Range("D20").Value = 6000
$ws.Range("E20").Value = -8900
$ws.Range("D21").Value = 2384500
$ws.Range("E21").Value = 2783900
$ws.Range("D22").Value = 254100
$ws.Range("E22").Value = 250100
$ws.Range("D23").Value = 1751400
$ws.Range("E23").Value = 2158700
$ws.Range("D24").Value = 360100
$ws.Range("E24").Value = 474400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1391300
$ws.Range("E26").Value = 1684300
$ws.Range("D27").Value = 1391300
$ws.Range("E27").Value = 1684300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 24000
$ws.Range("E29").Value = -47000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -6000
$ws.Range("E32").Value = 8900
$ws.Range("D33").Value = 1415300
$ws.Range("E33").Value = 1637300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 1415300
$ws.Range("E35").Value = 1637300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 866000
$ws.Range("E41").Value = 2574500
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 2441500
$ws.Range("E43").Value = 2266800
$ws.Range("D44").Value = 51100
$ws.Range("E44").Value = 41900
$ws.Range("D45").Value = 694600
$ws.Range("E45").Value = 669900
$ws.Range("D46").Value = 4053200
$ws.Range("E46").Value = 5553100
$ws.Range("D47").Value = 1202800
$ws.Range("E47").Value = 1135500
$ws.Range("D48").Value = 22842700
$ws.Range("E48").Value = 22613600
$ws.Range("D49").Value = 2331500
$ws.Range("E49").Value = 2345000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 2381000
$ws.Range("E52").Value = 2406500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 32811200
$ws.Range("E54").Value = 34053700
$ws.Range("D57").Value = 1207900
$ws.Range("E57").Value = 932800
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 1765600
$ws.Range("E59").Value = 2693900
$ws.Range("D60").Value = 2973500
$ws.Range("E60").Value = 3626700
$ws.Range("D61").Value = 31075300
$ws.Range("E61").Value = 31895200
$ws.Range("D62").Value = 5020800
$ws.Range("E62").Value = 5324400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 39069600
$ws.Range("E66").Value = 40846300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 50487000
$ws.Range("E72").Value = 49076200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -6258400
$ws.Range("E76").Value = -6792600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 1415300
$ws.Range("E81").Value = 1637300
$ws.Range("D83").Value = 379000
$ws.Range("E83").Value = 375100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 1511600
$ws.Range("E89").Value = 2471100
$ws.Range("D91").Value = -873800
$ws.Range("E91").Value = -703800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -883900
$ws.Range("E94").Value = -721300
$ws.Range("D96").Value = -892500
$ws.Range("E96").Value = -779800
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -2311200
$ws.Range("E100").Value = -768000
$ws.Range("D101").Value = -25000
$ws.Range("E101").Value = -30800
$ws.Range("D102").Value = -1708500
$ws.Range("E102").Value = 951000

